$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.440.49"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.823.77"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4577"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3811"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07880"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9682"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.78%  "
$ws.Range("D13").Value = "1.834.93"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.877"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.030"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06616"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001026"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").Value = "27.436.94"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.332"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.300"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").Value = "2.047.76"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.054"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.263"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9446"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09315"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.237"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.321"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05918"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02185"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.156"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.001"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5754"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1829"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.263"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5438"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.865"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06609"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.040"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.37%  "
